# Update 27 mai 2025
# - Add vertical-center alignment to the whole "JUDUL" (F) column so it is
#   consistent with the other columns (D, E, G, H, I already had it).
# - Fill in thesis data that had become available for NIM 19065008 (date),
#   NIM 19065014 and NIM 19065028 (title / pembimbing / penelaah / date).
# - Highlight NIM 19065022 in yellow, like the other "not submitted yet" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply vertical-center alignment across column F (JUDUL) -----------
# This matches the rest of the table (columns D, E, G, H, I already use
# vertical-center) and is applied at the column level so both filled and
# still-empty cells pick it up.
$ws.Columns("F").VerticalAlignment = -4108   # xlCenter

# --- 2. NIM 19065008 (Igo Firnandes) now has a graduation date ------------
$ws.Range("I10").Value = "12/1/2024"
$ws.Range("I10").NumberFormat = "mmm-yy"

# --- 3. NIM 19065014 (Putri Ardini) thesis info filled in ------------------
$ws.Range("F16").Value = "HUBUNGAN KESIAPAN GURU DAN SARANA PEMBELAJARAN TERHADAP IMPLEMENTASI KURIKULUM MARDEKA DI JURUSAN TEKNIK ELEKTRONIKA PADA SMK NEGERI 5 PADANG"
$ws.Range("G16").Value = "Drs. Efrizon, MT / 5319"
$ws.Range("H16").Value = "Dr. Dedy Irfan, S.Pd., M.Kom / 5327, Ilmiyati Rahmy Jasril, S.Pd., M.Pd.T / 5340"
$ws.Range("I16").Value = "3/1/2025"
$ws.Range("I16").NumberFormat = "mmm-yy"
$ws.Range("F16:I16").WrapText = $true
$ws.Range("F16:I16").VerticalAlignment = -4108   # xlCenter
$ws.Rows(16).RowHeight = 45

# --- 4. NIM 19065028 (ARIZAL ALFANDI) thesis info filled in ----------------
$ws.Range("F30").Value = "Prototipe mesin pemotong rumput berbasis arduino uno dengan pengendali android"
$ws.Range("G30").Value = "Thamrin S.Pd., M.T. / 5334"
$ws.Range("H30").Value = "Sartika Anori, S.Pd.,M.Pd.T / 182038, Winda Agustiarmi, S.Pd.,M.Pd.T / 192042"
$ws.Range("F30:H30").WrapText = $true
$ws.Range("F30:H30").VerticalAlignment = -4108   # xlCenter
$ws.Rows(30).RowHeight = 45

# --- 5. Highlight NIM 19065022 (Serli Rissandi) in yellow, like the other
#        rows that are still awaiting thesis data -----------------------
$ws.Range("D24:I24").Interior.Color = 65535   # yellow, matches rows 4/20/27

# --- 6. Update the view so it reflects where the editor left off ----------
$ws.Range("F29").Select()
